{"js": "// Replace each known old cell/date string with its new value.\n// Uses Word.js range.search() + insertText(..., \"Replace\") for exact,\n// unambiguous text substitutions (each old string occurs exactly once\n// in the document, so a direct search+replace is unambiguous).\nconst replacements = [\n  [\"2024-04-07 Sunday\", \"2024-04-08 Monday\"],\n  [\"315\u00f74=78, 3\", \"950\u00f77=135, 5\"],\n  [\"189\u00f77=27, 0\", \"154\u00f73=51, 1\"],\n  [\"151\u00f75=30, 1\", \"721\u00f79=80, 1\"],\n  [\"128\u00f76=21, 2\", \"721\u00f78=90, 1\"],\n  [\"211\u00f72=105, 1\", \"635\u00f75=127, 0\"],\n  [\"958\u00f78=119, 6\", \"270\u00f75=54, 0\"],\n  [\"513\u00f79=57, 0\", \"227\u00f76=37, 5\"],\n  [\"260\u00f74=65, 0\", \"388\u00f72=194, 0\"],\n  [\"520\u00f74=130, 0\", \"925\u00f75=185, 0\"],\n  [\"119\u00f73=39, 2\", \"326\u00f74=81, 2\"],\n  [\"929\u00f77=132, 5\", \"180\u00f72=90, 0\"],\n  [\"825\u00f76=137, 3\", \"169\u00f72=84, 1\"],\n  [\"698\u00f73=232, 2\", \"404\u00f76=67, 2\"],\n  [\"485\u00f72=242, 1\", \"940\u00f73=313, 1\"],\n  [\"132\u00f76=22, 0\", \"397\u00f75=79, 2\"],\n  [\"547\u00f72=273, 1\", \"480\u00f78=60, 0\"],\n  [\"100\u00f78=12, 4\", \"780\u00f74=195, 0\"],\n  [\"753\u00f78=94, 1\", \"487\u00f77=69, 4\"],\n  [\"134\u00f79=14, 8\", \"131\u00f75=26, 1\"],\n  [\"912\u00f79=101, 3\", \"226\u00f73=75, 1\"],\n  [\"808\u00f73=269, 1\", \"929\u00f74=232, 1\"],\n  [\"707\u00f77=101, 0\", \"954\u00f79=106, 0\"],\n  [\"745\u00f72=372, 1\", \"553\u00f73=184, 1\"],\n  [\"489\u00f79=54, 3\", \"872\u00f75=174, 2\"],\n  [\"124\u00f74=31, 0\", \"126\u00f77=18, 0\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each known old cell/date string with its new value using\n# Word's Find/Replace COM object model (Range.Find / Find.Execute).\n# Each old string is unique in the document, so wdReplaceAll (2) on a\n# Find scoped to the whole document body is an exact, unambiguous swap.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"2024-04-07 Sunday\"; New = \"2024-04-08 Monday\" },\n  @{ Old = \"315\u00f74=78, 3\"; New = \"950\u00f77=135, 5\" },\n  @{ Old = \"189\u00f77=27, 0\"; New = \"154\u00f73=51, 1\" },\n  @{ Old = \"151\u00f75=30, 1\"; New = \"721\u00f79=80, 1\" },\n  @{ Old = \"128\u00f76=21, 2\"; New = \"721\u00f78=90, 1\" },\n  @{ Old = \"211\u00f72=105, 1\"; New = \"635\u00f75=127, 0\" },\n  @{ Old = \"958\u00f78=119, 6\"; New = \"270\u00f75=54, 0\" },\n  @{ Old = \"513\u00f79=57, 0\"; New = \"227\u00f76=37, 5\" },\n  @{ Old = \"260\u00f74=65, 0\"; New = \"388\u00f72=194, 0\" },\n  @{ Old = \"520\u00f74=130, 0\"; New = \"925\u00f75=185, 0\" },\n  @{ Old = \"119\u00f73=39, 2\"; New = \"326\u00f74=81, 2\" },\n  @{ Old = \"929\u00f77=132, 5\"; New = \"180\u00f72=90, 0\" },\n  @{ Old = \"825\u00f76=137, 3\"; New = \"169\u00f72=84, 1\" },\n  @{ Old = \"698\u00f73=232, 2\"; New = \"404\u00f76=67, 2\" },\n  @{ Old = \"485\u00f72=242, 1\"; New = \"940\u00f73=313, 1\" },\n  @{ Old = \"132\u00f76=22, 0\"; New = \"397\u00f75=79, 2\" },\n  @{ Old = \"547\u00f72=273, 1\"; New = \"480\u00f78=60, 0\" },\n  @{ Old = \"100\u00f78=12, 4\"; New = \"780\u00f74=195, 0\" },\n  @{ Old = \"753\u00f78=94, 1\"; New = \"487\u00f77=69, 4\" },\n  @{ Old = \"134\u00f79=14, 8\"; New = \"131\u00f75=26, 1\" },\n  @{ Old = \"912\u00f79=101, 3\"; New = \"226\u00f73=75, 1\" },\n  @{ Old = \"808\u00f73=269, 1\"; New = \"929\u00f74=232, 1\" },\n  @{ Old = \"707\u00f77=101, 0\"; New = \"954\u00f79=106, 0\" },\n  @{ Old = \"745\u00f72=372, 1\"; New = \"553\u00f73=184, 1\" },\n  @{ Old = \"489\u00f79=54, 3\"; New = \"872\u00f75=174, 2\" },\n  @{ Old = \"124\u00f74=31, 0\"; New = \"126\u00f77=18, 0\" },\n)\n\nforeach ($item in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $item.Old\n  $find.Replacement.Text = $item.New\n  $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n  if (-not $found) {\n    throw \"Could not find text: \" + $item.Old\n  }\n}\n"}
